# Update cryptos list with refreshed prices / 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2..44 (excluding the two swapped pairs handled below) ---
$ws.Range("D2").Value = "69.738.95"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "2.514.51"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "576.50"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").Value = "166.77"
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "2.511.45"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +3.36%  "

$ws.Range("E13").Value = "  +2.49%  "

$ws.Range("D14").Value = "2.977.83"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").Value = "69.590.71"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("D17").Value = "24.99"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "2.513.93"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("E20").Value = "  +4.60%  "

$ws.Range("D21").Value = "350.46"
$ws.Range("E21").Value = "  -2.84%  "

$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").Value = "2.01"
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "70.31"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("D26").Value = "4.02"
$ws.Range("E26").Value = "  -1.70%  "

$ws.Range("D27").Value = "8.95"
$ws.Range("E27").Value = "  -3.30%  "

$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").Value = "0.0₃0901"
$ws.Range("E30").Value = "  -2.06%  "

$ws.Range("D31").Value = "7.94"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").Value = "466.64"
$ws.Range("E32").Value = "  -3.77%  "

$ws.Range("D33").Value = "1.26"
$ws.Range("E33").Value = "  -1.68%  "

$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").Value = "157.94"
$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").Value = "19.02"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").Value = "18.57"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").Value = "4.78"
$ws.Range("E41").Value = "  +0.51%  "

$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("E43").Value = "  -3.68%  "

$ws.Range("D44").Value = "38.41"
$ws.Range("E44").Value = "  +0.03%  "

# --- Rows 45/46: ImmutableX and dogwifhat swap places (ranking reorder) ---
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -6.18%  "

$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  -13.18%  "

# --- Row 47..49 ---
$ws.Range("D47").Value = "142.60"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").Value = "0.529"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("E49").Value = "  -1.00%  "

# --- Rows 50/51: Cronos and Optimism swap places (ranking reorder) ---
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.57"
$ws.Range("E50").Value = "  -3.30%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0729"
$ws.Range("E51").Value = "  -0.84%  "
